# Auto-generated edits applying Ultima_Profits market-price refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H5").Value = 343.75
$ws.Range("I5").Value = 128.875
$ws.Range("J5").Value = 773.5
$ws.Range("K5").Value = 128.875
$ws.Range("L5").Value = 773.5
$ws.Range("M5").Value = -13.875
$ws.Range("N5").Value = -1003.5

$ws.Range("H64").Value = 2979.16
$ws.Range("I64").Value = 2810.2
$ws.Range("J64").Value = 3091.8
$ws.Range("K64").Value = 2810.2
$ws.Range("L64").Value = 3091.8
$ws.Range("M64").Value = -2562.2
$ws.Range("N64").Value = -3587.8

$ws.Range("H67").Value = 2979.16
$ws.Range("I67").Value = 2810.2
$ws.Range("J67").Value = 3091.8
$ws.Range("K67").Value = 2810.2
$ws.Range("L67").Value = 3091.8
$ws.Range("M67").Value = -1952.2
$ws.Range("N67").Value = -4807.8

$ws.Range("H80").Value = 1703.4667
$ws.Range("I80").Value = 2183.1667
$ws.Range("J80").Value = 1383.6666
$ws.Range("K80").Value = 6549.500100000001
$ws.Range("L80").Value = 4150.9998
$ws.Range("M80").Value = -5551.500100000001
$ws.Range("N80").Value = -6146.9998

$ws.Range("H83").Value = 1703.4667
$ws.Range("I83").Value = 2183.1667
$ws.Range("J83").Value = 1383.6666
$ws.Range("K83").Value = 19648.5003
$ws.Range("L83").Value = 12452.9994
$ws.Range("M83").Value = -14656.5003
$ws.Range("N83").Value = -22436.9994

$ws.Range("H103").Value = 4623615
$ws.Range("I103").Value = 12020401
$ws.Range("J103").Value = 623.75
$ws.Range("K103").Value = 36061203
$ws.Range("L103").Value = 1871.25
$ws.Range("M103").Value = -36060617
$ws.Range("N103").Value = -3043.25

$ws.Range("H111").Value = 2189.4443
$ws.Range("I111").Value = 3829.6667
$ws.Range("J111").Value = 1369.3334
$ws.Range("K111").Value = 11489.0001
$ws.Range("L111").Value = 4108.0002
$ws.Range("M111").Value = -8422.000100000001
$ws.Range("N111").Value = -10242.0002

$ws.Range("H127").Value = 538563.9
$ws.Range("J127").Value = 673073.6
$ws.Range("L127").Value = 2019220.8
$ws.Range("N127").Value = -2029140.8

$ws.Range("H129").Value = 1119.84
$ws.Range("I129").Value = 643.5
$ws.Range("J129").Value = 1139.6875
$ws.Range("K129").Value = 1930.5
$ws.Range("L129").Value = 3419.0625
$ws.Range("M129").Value = 3069.5
$ws.Range("N129").Value = -13419.0625

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 1888.625
$ws.Range("I2").Value = 1905.5
$ws.Range("J2").Value = 1883
$ws.Range("K2").Value = 1905.5
$ws.Range("L2").Value = 1883
$ws.Range("M2").Value = -1792.5
$ws.Range("N2").Value = -2109

$ws.Range("H97").Value = 9653.923000000001
$ws.Range("I97").Value = 13032.223
$ws.Range("J97").Value = 2052.75
$ws.Range("K97").Value = 13032.223
$ws.Range("L97").Value = 2052.75
$ws.Range("M97").Value = -12536.223
$ws.Range("N97").Value = -3044.75

$ws.Range("H102").Value = 5042.484
$ws.Range("I102").Value = 5848.524
$ws.Range("J102").Value = 3349.8
$ws.Range("K102").Value = 5848.524
$ws.Range("L102").Value = 3349.8
$ws.Range("M102").Value = -4226.524
$ws.Range("N102").Value = -6593.8

$ws.Range("H116").Value = 1888.625
$ws.Range("I116").Value = 1905.5
$ws.Range("J116").Value = 1883
$ws.Range("K116").Value = 1905.5
$ws.Range("L116").Value = 1883
$ws.Range("M116").Value = 388.5
$ws.Range("N116").Value = -6471

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 1888.625
$ws.Range("I3").Value = 1905.5
$ws.Range("J3").Value = 1883
$ws.Range("K3").Value = 1905.5
$ws.Range("L3").Value = 1883
$ws.Range("M3").Value = -1791.5
$ws.Range("N3").Value = -2111

$ws.Range("H7").Value = 776525.9
$ws.Range("I7").Value = 1018533.8
$ws.Range("J7").Value = 50502
$ws.Range("K7").Value = 1018533.8
$ws.Range("L7").Value = 50502
$ws.Range("M7").Value = -1018420.8
$ws.Range("N7").Value = -50728

$ws.Range("H99").Value = 1795.8823
$ws.Range("I99").Value = 1877.5
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 1877.5
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -379.5
$ws.Range("N99").Value = -4596

$ws.Range("H105").Value = 4406.875
$ws.Range("I105").Value = 2949.5
$ws.Range("J105").Value = 4615.0713
$ws.Range("K105").Value = 2949.5
$ws.Range("L105").Value = 4615.0713
$ws.Range("M105").Value = -1202.5
$ws.Range("N105").Value = -8109.0713

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 988.8889
$ws.Range("I16").Value = 650
$ws.Range("J16").Value = 1666.6666
$ws.Range("K16").Value = 650
$ws.Range("L16").Value = 1666.6666
$ws.Range("M16").Value = -363
$ws.Range("N16").Value = -2240.6666

$ws.Range("H31").Value = 7095921
$ws.Range("I31").Value = 6294.5654
$ws.Range("J31").Value = 13890146
$ws.Range("K31").Value = 6294.5654
$ws.Range("L31").Value = 13890146
$ws.Range("M31").Value = -5999.5654
$ws.Range("N31").Value = -13890736

$ws.Range("H34").Value = 7095921
$ws.Range("I34").Value = 6294.5654
$ws.Range("J34").Value = 13890146
$ws.Range("K34").Value = 6294.5654
$ws.Range("L34").Value = 13890146
$ws.Range("M34").Value = -6092.5654
$ws.Range("N34").Value = -13890550

$ws.Range("H58").Value = 3256.9443
$ws.Range("I58").Value = 1894.5
$ws.Range("K58").Value = 1894.5
$ws.Range("M58").Value = -1691.5

$ws.Range("H94").Value = 3209.3076
$ws.Range("I94").Value = 1841.6428
$ws.Range("J94").Value = 4804.9165
$ws.Range("K94").Value = 1841.6428
$ws.Range("L94").Value = 4804.9165
$ws.Range("M94").Value = -1390.6428
$ws.Range("N94").Value = -5706.9165

$ws.Range("H99").Value = 1163.75
$ws.Range("I99").Value = 1163.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1163.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 334.25
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 864.875
$ws.Range("I105").Value = 845.5714
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 845.5714
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 901.4286
$ws.Range("N105").Value = -4494

$ws.Range("H107").Value = 497.5
$ws.Range("I107").Value = 307.4
$ws.Range("J107").Value = 756.7273
$ws.Range("K107").Value = 307.4
$ws.Range("L107").Value = 756.7273
$ws.Range("M107").Value = 1612.6
$ws.Range("N107").Value = -4596.7273

$ws.Range("H113").Value = 988.8889
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 1666.6666
$ws.Range("K113").Value = 650
$ws.Range("L113").Value = 1666.6666
$ws.Range("M113").Value = 1520
$ws.Range("N113").Value = -6006.6666

$ws.Range("H126").Value = 1163.75
$ws.Range("I126").Value = 1163.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3491.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1021.25
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 3256.9443
$ws.Range("I136").Value = 1894.5
$ws.Range("K136").Value = 5683.5
$ws.Range("M136").Value = -3133.5

$ws.Range("H138").Value = 74221
$ws.Range("J138").Value = 79748.625
$ws.Range("L138").Value = 79748.625
$ws.Range("N138").Value = -90028.625

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H49").Value = 5026.533
$ws.Range("J49").Value = 5026.533
$ws.Range("L49").Value = 15079.599
$ws.Range("N49").Value = -15391.599

$ws.Range("H68").Value = 1198.1818
$ws.Range("I68").Value = 964.1579
$ws.Range("J68").Value = 1376.04
$ws.Range("K68").Value = 2892.4737
$ws.Range("L68").Value = 4128.12
$ws.Range("M68").Value = -2081.4737
$ws.Range("N68").Value = -5750.12

$ws.Range("H71").Value = 1198.1818
$ws.Range("I71").Value = 964.1579
$ws.Range("J71").Value = 1376.04
$ws.Range("K71").Value = 8677.4211
$ws.Range("L71").Value = 12384.36
$ws.Range("M71").Value = -4621.4211
$ws.Range("N71").Value = -20496.36

$ws.Range("H131").Value = 732.85
$ws.Range("J131").Value = 813.0741
$ws.Range("L131").Value = 2439.2223
$ws.Range("N131").Value = -12519.2223

$ws.Range("H139").Value = 1073.75
$ws.Range("I139").Value = 1073.75
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3221.25
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 1918.75
$ws.Range("N139").ClearContents()

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 6477.6665
$ws.Range("I102").Value = 7585.5713
$ws.Range("J102").Value = 2600
$ws.Range("K102").Value = 7585.5713
$ws.Range("L102").Value = 2600
$ws.Range("M102").Value = -5963.5713
$ws.Range("N102").Value = -5844

$ws.Range("H126").Value = 4699.8887
$ws.Range("I126").Value = 2925
$ws.Range("J126").Value = 5207
$ws.Range("K126").Value = 8775
$ws.Range("L126").Value = 15621
$ws.Range("M126").Value = -6305
$ws.Range("N126").Value = -20561

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H40").Value = 11998.75
$ws.Range("I40").Value = 14747.5
$ws.Range("K40").Value = 14747.5
$ws.Range("M40").Value = -14611.5

$ws.Range("H46").Value = 491.92307
$ws.Range("I46").Value = 528.4167
$ws.Range("K46").Value = 528.4167
$ws.Range("M46").Value = -340.4167

$ws.Range("H100").Value = 1527.4706
$ws.Range("I100").Value = 1464.7778
$ws.Range("J100").Value = 1598
$ws.Range("K100").Value = 1464.7778
$ws.Range("L100").Value = 1598
$ws.Range("M100").Value = -923.7778000000001
$ws.Range("N100").Value = -2680

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H21").Value = 12000017
$ws.Range("J21").Value = 12000017
$ws.Range("L21").Value = 12000017
$ws.Range("N21").Value = -12000487

$ws.Range("H24").Value = 12000010
$ws.Range("J24").Value = 12000010
$ws.Range("L24").Value = 12000010
$ws.Range("N24").Value = -12000470

$ws.Range("H25").Value = 9168360
$ws.Range("J25").Value = 9168360
$ws.Range("L25").Value = 9168360
$ws.Range("N25").Value = -9168946

$ws.Range("H28").Value = 4055.2
$ws.Range("I28").Value = 200
$ws.Range("J28").Value = 5019
$ws.Range("K28").Value = 200
$ws.Range("L28").Value = 5019
$ws.Range("M28").Value = 148
$ws.Range("N28").Value = -5715

$ws.Range("H35").Value = 12000017
$ws.Range("J35").Value = 12000017
$ws.Range("L35").Value = 12000017
$ws.Range("N35").Value = -12000597

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H137").Value = 64681.8
$ws.Range("J137").Value = 64681.8
$ws.Range("L137").Value = 64681.8
$ws.Range("N137").Value = -74881.8
